$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells I1 and J1, matching the style of existing header cell H1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Add the new data columns I and J for rows 2-12
$data = @(
    @(1, 4),
    @(1, 5),
    @(1, 3),
    @(1, 6),
    @(1, 5),
    @(1, 5),
    @(1, 5),
    @(1, 6),
    @(1, 6),
    @(1, 4),
    @(6, 8)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
